$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$url = "https://en.wikipedia.org/w/index.php?title=Anopheles_arabiensis&action=edit&redlink=1"
$ws.Hyperlinks.Add($ws.Range("Z1"), $url, "", "tooltip-stage", $url)
$ws.Range("Z1").Copy($ws.Range("A3"))
